$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$tbl.ApplyStyle("{8A4F1401-42BC-4CDB-8861-DD05B5FF2E5F}")
